# The "Date" / ${SALE_DATE} column is removed from the mandate table and
# the remaining five columns are widened to fill the freed space.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Drop the first column (header "Date", data "${SALE_DATE}").
$t.Columns.Item(1).Delete()

# Resize the table and the remaining columns to the new target widths
# (values are in dxa/twips in the XML; Word's COM Width/PreferredWidth
# properties are expressed in points, so divide by 20).
$t.PreferredWidthType = 3
$t.PreferredWidth = 9353 / 20

$widths = @(1870, 1871, 1870, 1871, 1871)
for ($i = 1; $i -le $t.Columns.Count; $i++) {
  $col = $t.Columns.Item($i)
  $col.PreferredWidthType = 3
  $col.Width = $widths[$i - 1] / 20
}
